$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.436.62'
$ws.Range("E2").Value = '  -3.14%  '
$ws.Range("D3").Value = '2.967.63'
$ws.Range("E3").Value = '  -5.42%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''494.99'
$ws.Range("E5").Value = '  -5.42%  '
$ws.Range("D6").Value = '''134.55'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '2.965.80'
$ws.Range("E8").Value = '  -5.55%  '
$ws.Range("D9").Value = '''0.423'
$ws.Range("E9").Value = '  -4.21%  '
$ws.Range("D10").Value = '''7.25'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").Value = '''0.103'
$ws.Range("E11").Value = '  -5.17%  '
$ws.Range("D12").Value = '''0.349'
$ws.Range("E12").Value = '  -8.70%  '
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '3.481.37'
$ws.Range("E14").Value = '  -5.18%  '
$ws.Range("D15").Value = '''24.78'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").Value = '56.457.66'
$ws.Range("E16").Value = '  -2.94%  '
$ws.Range("D17").Value = '2.978.42'
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").Value = '''0.0000145'
$ws.Range("E18").Value = '  -5.06%  '
$ws.Range("D19").Value = '''5.79'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '''12.18'
$ws.Range("E20").Value = '  -6.55%  '
$ws.Range("D21").Value = '''7.67'
$ws.Range("E21").Value = '  -3.85%  '
$ws.Range("D22").Value = '''321.51'
$ws.Range("E22").Value = '  -6.91%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").Value = '''0.459'
$ws.Range("E24").Value = '  -9.47%  '
$ws.Range("D25").Value = '''61.16'
$ws.Range("E25").Value = '  -11.14%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("D28").Value = '0.0₃0890'
$ws.Range("E28").Value = '  -7.29%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''6.46'
$ws.Range("E30").Value = '  -4.52%  '
$ws.Range("D31").Value = '''6.73'
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("D32").Value = '''1.17'
$ws.Range("E32").Value = '  -5.02%  '
$ws.Range("D33").Value = '''1.72'
$ws.Range("E33").Value = '  -7.25%  '
$ws.Range("D34").Value = '''19.75'
$ws.Range("E34").Value = '  -8.41%  '
$ws.Range("D35").Value = '''154.22'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("D36").Value = '''4.43'
$ws.Range("E36").Value = '  -7.39%  '
$ws.Range("D37").Value = '''1.27'
$ws.Range("E37").Value = '  -6.84%  '
$ws.Range("D38").Value = '''5.58'
$ws.Range("E38").Value = '  -10.14%  '
$ws.Range("D39").Value = '''0.0665'
$ws.Range("E39").Value = '  -4.04%  '
$ws.Range("D40").Value = '''23.10'
$ws.Range("E40").Value = '  -4.35%  '
$ws.Range("D41").Value = '3.000.20'
$ws.Range("E41").Value = '  -5.31%  '
$ws.Range("D42").Value = '''36.96'
$ws.Range("E42").Value = '  -8.67%  '
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").Value = '''0.994'
$ws.Range("E44").Value = '  -7.82%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '''0.633'
$ws.Range("E45").Value = '  -9.37%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '''1.40'
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("D47").Value = '2.205.73'
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").Value = '''3.54'
$ws.Range("E48").Value = '  -9.47%  '
$ws.Range("D49").Value = '''1.92'
$ws.Range("E49").Value = '  +5.13%  '
$ws.Range("D50").Value = '''0.0235'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '''5.62'
$ws.Range("E51").Value = '  -9.14%  '
